$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 562731.75
$ws.Range("I88").Value = 455.5
$ws.Range("J88").Value = 1237463.2
$ws.Range("K88").Value = 455.5
$ws.Range("L88").Value = 1237463.2
$ws.Range("M88").Value = -49.5
$ws.Range("N88").Value = -1238275.2
$ws.Range("H91").Value = 562731.75
$ws.Range("I91").Value = 455.5
$ws.Range("J91").Value = 1237463.2
$ws.Range("K91").Value = 455.5
$ws.Range("L91").Value = 1237463.2
$ws.Range("M91").Value = 948.5
$ws.Range("N91").Value = -1240271.2
$ws.Range("H103").Value = 1481.4286
$ws.Range("J103").Value = 2766.6667
$ws.Range("L103").Value = 8300.000100000001
$ws.Range("N103").Value = -9472.000100000001
$ws.Range("H107").Value = 1832.2059
$ws.Range("I107").Value = 1024.9286
$ws.Range("J107").Value = 5599.5
$ws.Range("K107").Value = 1024.9286
$ws.Range("L107").Value = 5599.5
$ws.Range("M107").Value = 895.0714
$ws.Range("N107").Value = -9439.5
$ws.Range("H112").Value = 2291.8333
$ws.Range("J112").Value = 2490.875
$ws.Range("L112").Value = 7472.625
$ws.Range("N112").Value = -9688.625
$ws.Range("H129").Value = 649.3182
$ws.Range("J129").Value = 919.5833
$ws.Range("L129").Value = 2758.7499
$ws.Range("N129").Value = -12758.7499
$ws.Range("H132").Value = 10108488
$ws.Range("I132").Value = 12826443
$ws.Range("K132").Value = 38479329
$ws.Range("M132").Value = -38476799
$ws.Range("H136").Value = 36662.855
$ws.Range("J136").Value = 36662.855
$ws.Range("L136").Value = 36662.855
$ws.Range("N136").Value = -46862.855
$ws.Range("H138").Value = 2071.2334
$ws.Range("J138").Value = 2251.1807
$ws.Range("L138").Value = 6753.5421
$ws.Range("N138").Value = -17033.5421
$ws.Range("H139").Value = 57170
$ws.Range("I139").Value = 20000
$ws.Range("K139").Value = 20000
$ws.Range("M139").Value = -14860
$ws.Range("H140").Value = 35987.777
$ws.Range("J140").Value = 35987.777
$ws.Range("L140").Value = 35987.777
$ws.Range("N140").Value = -46347.777

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1046.579
$ws.Range("I2").Value = 752.3333
$ws.Range("J2").Value = 2150
$ws.Range("K2").Value = 752.3333
$ws.Range("L2").Value = 2150
$ws.Range("M2").Value = -639.3333
$ws.Range("N2").Value = -2376
$ws.Range("H24").Value = 23538.75
$ws.Range("J24").Value = 23538.75
$ws.Range("L24").Value = 23538.75
$ws.Range("N24").Value = -24286.75
$ws.Range("H45").Value = 1274.8572
$ws.Range("I45").Value = 1144.8
$ws.Range("K45").Value = 1144.8
$ws.Range("M45").Value = -767.8
$ws.Range("H74").Value = 1814.5
$ws.Range("I74").Value = 1173.2778
$ws.Range("K74").Value = 1173.2778
$ws.Range("M74").Value = -299.2778000000001
$ws.Range("H77").Value = 1814.5
$ws.Range("I77").Value = 1173.2778
$ws.Range("K77").Value = 5866.389
$ws.Range("M77").Value = -1498.389
$ws.Range("H92").Value = 1267775
$ws.Range("J92").Value = 1267775
$ws.Range("L92").Value = 1267775
$ws.Range("N92").Value = -1272767
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").Value = $null
$ws.Range("H96").Value = 20000
$ws.Range("J96").Value = 20000
$ws.Range("L96").Value = 20000
$ws.Range("N96").Value = -25492
$ws.Range("H97").Value = 717.1429000000001
$ws.Range("I97").Value = 717.1429000000001
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 717.1429000000001
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -221.1429000000001
$ws.Range("N97").Value = $null
$ws.Range("H100").Value = 23538.75
$ws.Range("J100").Value = 23538.75
$ws.Range("L100").Value = 23538.75
$ws.Range("N100").Value = -25702.75
$ws.Range("H101").Value = 37000
$ws.Range("J101").Value = 37000
$ws.Range("L101").Value = 37000
$ws.Range("N101").Value = -43490
$ws.Range("H116").Value = 1046.579
$ws.Range("I116").Value = 752.3333
$ws.Range("J116").Value = 2150
$ws.Range("K116").Value = 752.3333
$ws.Range("L116").Value = 2150
$ws.Range("M116").Value = 1541.6667
$ws.Range("N116").Value = -6738
$ws.Range("H132").Value = 2638.8076
$ws.Range("I132").Value = 1939.3235
$ws.Range("K132").Value = 5817.970499999999
$ws.Range("M132").Value = -3287.970499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1046.579
$ws.Range("I3").Value = 752.3333
$ws.Range("J3").Value = 2150
$ws.Range("K3").Value = 752.3333
$ws.Range("L3").Value = 2150
$ws.Range("M3").Value = -638.3333
$ws.Range("N3").Value = -2378
$ws.Range("H107").Value = 877.9643
$ws.Range("I107").Value = 625.9091
$ws.Range("J107").Value = 1802.1666
$ws.Range("K107").Value = 625.9091
$ws.Range("L107").Value = 1802.1666
$ws.Range("M107").Value = 1294.0909
$ws.Range("N107").Value = -5642.1666
$ws.Range("H134").Value = 1025.9048
$ws.Range("I134").Value = 1047.2
$ws.Range("J134").Value = 600
$ws.Range("K134").Value = 3141.6
$ws.Range("L134").Value = 1800
$ws.Range("M134").Value = -606.6000000000004
$ws.Range("N134").Value = -6870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = $null
$ws.Range("H31").Value = 1291.1187
$ws.Range("I31").Value = 1181.4073
$ws.Range("K31").Value = 1181.4073
$ws.Range("M31").Value = -886.4073000000001
$ws.Range("H34").Value = 1291.1187
$ws.Range("I34").Value = 1181.4073
$ws.Range("K34").Value = 1181.4073
$ws.Range("M34").Value = -979.4073000000001
$ws.Range("H96").Value = 10397.5
$ws.Range("J96").Value = 10397.5
$ws.Range("L96").Value = 10397.5
$ws.Range("N96").Value = -15889.5
$ws.Range("H107").Value = 724.6
$ws.Range("I107").Value = 454.93332
$ws.Range("J107").Value = 1129.1
$ws.Range("K107").Value = 454.93332
$ws.Range("L107").Value = 1129.1
$ws.Range("M107").Value = 1465.06668
$ws.Range("N107").Value = -4969.1
$ws.Range("H132").Value = 1486.7805
$ws.Range("I132").Value = 1118.6666
$ws.Range("K132").Value = 3355.9998
$ws.Range("M132").Value = -825.9998000000001
$ws.Range("H134").Value = 17858696
$ws.Range("I134").Value = 1716.8334
$ws.Range("J134").Value = 50001260
$ws.Range("K134").Value = 5150.5002
$ws.Range("L134").Value = 150003780
$ws.Range("M134").Value = -2615.5002
$ws.Range("N134").Value = -150008850

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 689.96875
$ws.Range("I113").Value = 440.5
$ws.Range("J113").Value = 706.6
$ws.Range("K113").Value = 1321.5
$ws.Range("L113").Value = 2119.8
$ws.Range("M113").Value = 848.5
$ws.Range("N113").Value = -6459.8
$ws.Range("H131").Value = 25038250
$ws.Range("J131").Value = 54449.145
$ws.Range("L131").Value = 163347.435
$ws.Range("N131").Value = -173427.435

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1596.0714
$ws.Range("I102").Value = 1557.0952
$ws.Range("J102").Value = 1713
$ws.Range("K102").Value = 1557.0952
$ws.Range("L102").Value = 1713
$ws.Range("M102").Value = 64.90480000000002
$ws.Range("N102").Value = -4957
$ws.Range("H113").Value = 1123.4375
$ws.Range("I113").Value = 887.6
$ws.Range("J113").Value = 1516.5
$ws.Range("K113").Value = 887.6
$ws.Range("L113").Value = 1516.5
$ws.Range("M113").Value = 1282.4
$ws.Range("N113").Value = -5856.5
$ws.Range("H122").Value = 1609.6
$ws.Range("I122").Value = 1609.6
$ws.Range("K122").Value = 4828.799999999999
$ws.Range("M122").Value = -2378.799999999999
$ws.Range("H132").Value = 3215.5
$ws.Range("I132").Value = 2867.7144
$ws.Range("K132").Value = 8603.143199999999
$ws.Range("M132").Value = -6073.143199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").Value = $null
$ws.Range("H94").Value = 5000
$ws.Range("J94").Value = 5000
$ws.Range("L94").Value = 5000
$ws.Range("N94").Value = -6352

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").Value = $null
$ws.Range("H107").Value = 587.5172
$ws.Range("I107").Value = 370.41666
$ws.Range("J107").Value = 740.7646999999999
$ws.Range("K107").Value = 1111.24998
$ws.Range("L107").Value = 2222.2941
$ws.Range("M107").Value = 808.7500199999999
$ws.Range("N107").Value = -6062.2941
$ws.Range("H113").Value = 530.64703
$ws.Range("I113").Value = 364.76923
$ws.Range("J113").Value = 1069.75
$ws.Range("K113").Value = 1094.30769
$ws.Range("L113").Value = 3209.25
$ws.Range("M113").Value = 1075.69231
$ws.Range("N113").Value = -7549.25
$ws.Range("H122").Value = 32896310
$ws.Range("I122").Value = 36766264
$ws.Range("K122").Value = 110298792
$ws.Range("M122").Value = -110296342
